$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data rows 2-7 (row 1 is the header).
# Two new rows (8 and 9) are being appended with the same look & feel
# (borders/fill/alignment) as the existing data rows, so copy the
# formatting from the last data row (row 7) down into rows 8:9 first.
$ws.Range("A7:K7").Copy()
$ws.Range("A8:K9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 8: 2024-12-10 22:53:53
$ws.Cells.Item(8, 1).Value = "2024-12-10 22:53:53"
$ws.Cells.Item(8, 2).Value = 28.895
$ws.Cells.Item(8, 3).Value = 19.413
$ws.Cells.Item(8, 4).Value = 2.28
$ws.Cells.Item(8, 5).Value = 1.414
$ws.Cells.Item(8, 6).Value = 364124.18
$ws.Cells.Item(8, 7).Value = 366798.3
$ws.Cells.Item(8, 8).Value = 0.8108
$ws.Cells.Item(8, 9).Value = 0.7528
$ws.Cells.Item(8, 10).Value = 10.536
$ws.Cells.Item(8, 11).Value = 9.866

# Row 9: 2024-12-11 21:31:15
$ws.Cells.Item(9, 1).Value = "2024-12-11 21:31:15"
$ws.Cells.Item(9, 2).Value = 42.614
$ws.Cells.Item(9, 3).Value = 33.039
$ws.Cells.Item(9, 4).Value = 3.313
$ws.Cells.Item(9, 5).Value = 2.516
$ws.Cells.Item(9, 6).Value = 365226.46
$ws.Cells.Item(9, 7).Value = 365639.93
$ws.Cells.Item(9, 8).Value = 0.8898
$ws.Cells.Item(9, 9).Value = 0.8409
$ws.Cells.Item(9, 10).Value = 11.585
$ws.Cells.Item(9, 11).Value = 10.901

$null = $ws.Range("A1").Select()
